# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-10-28 21:13:17
#
# For every row in the "Recorded By" column (G), the comma-separated list of
# recorder names/emails that includes a "System" (or "system") entry gets its
# order reversed (e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ([string]::IsNullOrEmpty($val)) {
        continue
    }

    if ($val -notmatch ",") {
        continue
    }

    if ($val -notmatch "(?i)system") {
        continue
    }

    $parts = $val -split ","
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    $reversedParts = @()
    for ($i = $parts.Length - 1; $i -ge 0; $i--) {
        $reversedParts += $parts[$i]
    }
    $newVal = [string]::Join(", ", $reversedParts)

    $cell.Value2 = $newVal
}
